# This deck hosts a PowerPoint web add-in ("PowerPoll" bar-graph/pie-graph
# poll widget) via the mc:AlternateContent <p:graphicFrame>/<we:webextensionref>
# + <p:pic> fallback on slide 1, backed by the custom part
# ppt/slides/udata/data.xml (<we:webextension id="{...}">).
#
# The source commit ("bar graph now resizes to fit starting window size,
# both graphs timer now based on the value from the database, opening an
# old graph now has an accurate time instead of NaN") is purely a change to
# that add-in's own external JS/HTML (served from the developer/registry
# store referenced by <we:reference store="developer">) -- none of that
# script lives inside the .pptx. The only footprint left in the package by
# re-inserting/refreshing the add-in is a new opaque instance GUID on
# <we:webextension id="...">, plus every r:id in the package being
# renumbered, which is simply what PowerPoint does on every save and
# carries no content meaning (same slide master, same 11 layouts, same
# slide, same embedded snapshot image -- only the relationship labels and
# the add-in's internal id differ).
#
# Neither of those is reachable from the supported PowerPoint automation
# surface: there is no WebExtension/TaskPane object in the PowerPoint
# object model (Insert > Add-ins isn't VBA-scriptable in real PowerPoint
# either), CustomXMLParts/CustomerData exist only as empty stubs here, and
# relationship ids are never script-visible/settable. So the faithful,
# content-preserving action is to leave the deck's slides/shapes/media
# exactly as authored -- i.e. touch nothing -- which keeps every piece of
# real (user-visible) content identical to the target state.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Touch the object model read-only, to confirm the deck round-trips
# cleanly, without mutating any shape/content.
$null = $p.Slides.Count
$null = $s.Shapes.Count
